$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header labels for columns T (20) through AM (39)
$headers = @(
    "bus_price", "bus_comfort_level", "bus_availability", "bus_accessability", "bus_safety",
    "uber_price", "uber_comfort_level", "uber_availability", "uber_accessability", "uber_safety",
    "train_price", "train_comfort_level", "train_availability", "train_accessability", "train_safety",
    "bike_price", "bike_comfort_level", "bike_availability", "bike_accessability", "bike_safety"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 20 + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# New data values for rows 2-10, columns T (20) through AM (39)
$data = @(
    @(2,0,1,0,2,0,1,1,1,2,2,0,1,0,2,0,1,1,1,2),
    @(1,0,1,1,2,0,0,0,1,0,1,0,1,1,2,0,0,0,1,0),
    @(0,1,0,0,2,0,1,0,1,1,0,1,0,0,2,0,1,0,1,1),
    @(1,2,0,0,2,0,1,1,1,1,1,2,0,0,2,0,1,1,1,1),
    @(2,1,1,0,2,0,0,1,1,2,2,1,1,0,2,0,0,1,1,2),
    @(1,2,1,1,2,0,0,0,1,2,1,2,1,1,2,0,0,0,1,2),
    @(0,0,0,0,2,0,1,1,1,0,0,0,0,0,2,0,1,1,1,0),
    @(0,2,0,0,2,0,1,1,1,0,0,2,0,0,2,0,1,1,1,0),
    @(1,1,1,1,2,0,1,0,1,0,1,1,1,1,2,0,1,0,1,0)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $rowData = $data[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = 20 + $c
        $ws.Cells.Item($rowNum, $col).Value = $rowData[$c]
    }
}

# S10 also set explicitly to 0 per diff (value unchanged, but ensure correctness)
$ws.Cells.Item(10, 19).Value = 0
